$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 681.6667
$ws.Range("J19").Value = 666.6667
$ws.Range("L19").Value = 666.6667
$ws.Range("N19").Value = -1016.6667

$ws.Range("H43").Value = 1967.6364
$ws.Range("I43").Value = 2958.4
$ws.Range("J43").Value = 1142
$ws.Range("K43").Value = 2958.4
$ws.Range("L43").Value = 1142
$ws.Range("M43").Value = -2889.4
$ws.Range("N43").Value = -1280

$ws.Range("H116").Value = 18184246
$ws.Range("I116").Value = 66668500
$ws.Range("J116").Value = 2650
$ws.Range("K116").Value = 66668500
$ws.Range("L116").Value = 2650
$ws.Range("M116").Value = -66665058
$ws.Range("N116").Value = -9534

$ws.Range("H132").Value = 1962.7693
$ws.Range("I132").Value = 1903.7646
$ws.Range("J132").Value = 2074.2222
$ws.Range("K132").Value = 5711.293799999999
$ws.Range("L132").Value = 6222.6666
$ws.Range("M132").Value = -3181.293799999999
$ws.Range("N132").Value = -11282.6666

$ws.Range("H141").Value = 9183.846
$ws.Range("I141").Value = 1400.909
$ws.Range("J141").Value = 51990
$ws.Range("K141").Value = 4202.727000000001
$ws.Range("L141").Value = 155970
$ws.Range("M141").Value = 977.2729999999992
$ws.Range("N141").Value = -166330

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14904.216
$ws.Range("I32").Value = 18269.719
$ws.Range("K32").Value = 18269.719
$ws.Range("M32").Value = -17982.719

$ws.Range("H74").Value = 970.56525
$ws.Range("I74").Value = 746.4375
$ws.Range("J74").Value = 1482.8572
$ws.Range("K74").Value = 746.4375
$ws.Range("L74").Value = 1482.8572
$ws.Range("M74").Value = 127.5625
$ws.Range("N74").Value = -3230.8572

$ws.Range("H77").Value = 970.56525
$ws.Range("I77").Value = 746.4375
$ws.Range("J77").Value = 1482.8572
$ws.Range("K77").Value = 3732.1875
$ws.Range("L77").Value = 7414.286
$ws.Range("M77").Value = 635.8125
$ws.Range("N77").Value = -16150.286

$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 25000
$ws.Range("J20").Value = 25000
$ws.Range("L20").Value = 25000
$ws.Range("N20").Value = -25472

$ws.Range("H30").Value = 25000
$ws.Range("J30").Value = 25000
$ws.Range("L30").Value = 25000
$ws.Range("N30").Value = -25182

$ws.Range("H31").Value = 1709.95
$ws.Range("I31").Value = 1423.5294
$ws.Range("K31").Value = 1423.5294
$ws.Range("M31").Value = -1128.5294

$ws.Range("H34").Value = 1709.95
$ws.Range("I34").Value = 1423.5294
$ws.Range("K34").Value = 1423.5294
$ws.Range("M34").Value = -1221.5294

$ws.Range("H128").Value = 25000
$ws.Range("J128").Value = 25000
$ws.Range("L128").Value = 25000
$ws.Range("N128").Value = -34960

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1380.5555
$ws.Range("I5").Value = 1380.5555
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 4141.666499999999
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = $null
$ws.Range("N5").Value = -4029.666499999999

$ws.Range("H12").Value = 2760743
$ws.Range("I12").Value = 401
$ws.Range("J12").Value = 3220800
$ws.Range("K12").Value = 1203
$ws.Range("L12").Value = 9662400
$ws.Range("M12").Value = -1030
$ws.Range("N12").Value = -9662746

$ws.Range("H122").Value = 235
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = $null

$ws.Range("H135").Value = 1380.5555
$ws.Range("I135").Value = 1380.5555
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 12424.9995
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = $null
$ws.Range("N135").Value = -9889.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 28888
$ws.Range("J63").Value = 28888
$ws.Range("L63").Value = 28888
$ws.Range("N63").Value = -30260

$ws.Range("H66").Value = 28888
$ws.Range("J66").Value = 28888
$ws.Range("L66").Value = 86664
$ws.Range("N66").Value = -93528

$ws.Range("H132").Value = 1887.7307
$ws.Range("I132").Value = 1312.8572
$ws.Range("J132").Value = 4302.2
$ws.Range("K132").Value = 3938.5716
$ws.Range("L132").Value = 12906.6
$ws.Range("M132").Value = -1408.5716
$ws.Range("N132").Value = -17966.6

$ws.Range("H134").Value = 500326
$ws.Range("J134").Value = 500326
$ws.Range("L134").Value = 1500978
$ws.Range("N134").Value = -1506048

$ws.Range("H135").Value = 69887.5
$ws.Range("J135").Value = 69887.5
$ws.Range("L135").Value = 69887.5
$ws.Range("N135").Value = -80027.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5123
$ws.Range("I40").Value = 5254.143
$ws.Range("J40").Value = 4205
$ws.Range("K40").Value = 5254.143
$ws.Range("L40").Value = 4205
$ws.Range("M40").Value = -5118.143
$ws.Range("N40").Value = -4477

$ws.Range("H76").Value = 7166.5
$ws.Range("J76").Value = 7999.8
$ws.Range("L76").Value = 7999.8
$ws.Range("N76").Value = -8675.799999999999

$ws.Range("H79").Value = 7166.5
$ws.Range("J79").Value = 7999.8
$ws.Range("L79").Value = 7999.8
$ws.Range("N79").Value = -10339.8

$ws.Range("H132").Value = 5893.625
$ws.Range("I132").Value = 6527.75
$ws.Range("J132").Value = 4625.375
$ws.Range("K132").Value = 19583.25
$ws.Range("L132").Value = 13876.125
$ws.Range("M132").Value = -17053.25
$ws.Range("N132").Value = -18936.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 33889.715
$ws.Range("J64").Value = 33889.715
$ws.Range("L64").Value = 33889.715
$ws.Range("N64").Value = -34385.715

$ws.Range("H67").Value = 33889.715
$ws.Range("J67").Value = 33889.715
$ws.Range("L67").Value = 33889.715
$ws.Range("N67").Value = -35605.715

$ws.Range("H92").Value = 24966.666
$ws.Range("J92").Value = 24966.666
$ws.Range("L92").Value = 24966.666
$ws.Range("N92").Value = -29958.666

$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680

$ws.Range("H132").Value = 2234.5103
$ws.Range("I132").Value = 1335.9615
$ws.Range("K132").Value = 4007.8845
$ws.Range("M132").Value = -1477.8845
